$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter new projection rows in the same order the original author typed them,
# so that shared-string ids line up with the target workbook.
# Row 3: wintri (proj4 value typed before the crs name)
$ws.Cells.Item(3, 2).Value = "'+proj=wintri +datum=WGS84 +no_defs +over"
$ws.Cells.Item(3, 1).Value = "wintri"

# Row 4: robinson
$ws.Cells.Item(4, 1).Value = "robinson"
$ws.Cells.Item(4, 2).Value = "'+proj=robin +lat_0=0 +lon_0=0 +x0=0 +y0=0"

# Row 5: equirec
$ws.Cells.Item(5, 1).Value = "equirec"
$ws.Cells.Item(5, 2).Value = "'+proj=longlat +ellps=WGS84 +datum=WGS84 +no_defs"

# Row 6: gallpeters
$ws.Cells.Item(6, 1).Value = "gallpeters"
$ws.Cells.Item(6, 2).Value = "'+proj=cea +lon_0=0 +lat_ts=45 +x_0=0 +y_0=0 +ellps=WGS84 +units=m +no_defs"

# Row 7: hobodyer
$ws.Cells.Item(7, 1).Value = "hobodyer"
$ws.Cells.Item(7, 2).Value = "'+proj=cea +lat_ts=37.5"

# Row 8: goode
$ws.Cells.Item(8, 1).Value = "goode"
$ws.Cells.Item(8, 2).Value = "'+proj=igh"

# Widen column B to fit the longer strings
$ws.Columns.Item(2).ColumnWidth = 109

# Sort rows 2:8 by column A ascending (crs name), matching the workbook's sortState
$sortRange = $ws.Range("A2:B8")
$keyRange = $ws.Range("A2:A8")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Update selection to match the saved view state
$ws.Range("A2:B8").Select()
